# Upgrade left table: add a new "2023" column (K) to the Chiatura
# average-monthly-remuneration table, mirroring the formatting of the
# preceding "2022" column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/number formats/borders) of the existing
# last column (J, year 2022) into the new column (K) so the new
# "2023" column matches the look of the rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new column with the 2023 figures.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1541.9
$ws.Range("K5").Value = 931.9
$ws.Range("K6").Value = 1608.3
